# Insert a new weekly price record for "Terminal Hortofrutícola Agro Chillán - Mango"
# right before the existing row 50. This pushes the former rows 50-82 down to 51-83
# (dimension grows from A1:T82 to A1:T83) and fills the newly opened row 50 with the
# new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 50 (shifts rows 50..82 down to 51..83).
$ws.Rows.Item(50).Insert()

# Populate the new row 50 with the new data point.
$ws.Range("A50").Value = 7
$ws.Range("B50").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C50").Value = "Ñuble"
$ws.Range("D50").Value = 44846
$ws.Range("E50").Value = 16
$ws.Range("F50").Value = "Fruta"
$ws.Range("G50").Value = 100108
$ws.Range("H50").Value = "Tropicales y subtropicales"
$ws.Range("I50").Value = 100108002
$ws.Range("J50").Value = "Mango"
$ws.Range("K50").Value = "Sin especificar"
$ws.Range("L50").Value = "Primera"
$ws.Range("M50").Value = 120
$ws.Range("N50").Value = 7500
$ws.Range("O50").Value = 8000
$ws.Range("P50").Value = 7750
$ws.Range("Q50").Value = "$/bandeja 4 kilos"
$ws.Range("R50").Value = "Brasil"
$ws.Range("S50").Value = 1938
$ws.Range("T50").Value = 4
